$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.202277660369873
$ws.Range("B1").Value = 2.226264953613281
$ws.Range("D1").Value = 1.448760509490967
$ws.Range("E1").Value = 0.9057269096374512
